# Apply the "Top 10 Cities Least Debt" ranking update:
# - Remove the "spokane" row (it drops off the bottom-10 least-debt ranking)
# - Every row below it shifts up by one
# - A new entry for "aurora, Colorado" is appended as the new last row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top 10 Cities Least Debt")

# Full replacement data set for rows 2-11 (row 1 is the header, untouched)
$data = @(
    @("fresno",        "California", -705893000, -1301.780353046185),
    @("chula vista",   "California", -133323756, -483.933778584392),
    @("riverside",     "California", -91523000,  -290.5362936250452),
    @("madison",       "Wisconsin",  -78290849,  -290.2144019512991),
    @("plano",         "Texas",      -43321968,  -151.75717153176),
    @("wichita",       "Kansas",     -13760335,  -34.61310235016237),
    @("santa clarita", "California", -2504506,   -10.95205943702745),
    @("garland",       "Texas",      8188580,    33.28677525701115),
    @("irving",        "Texas",      12758086,   49.70289106966851),
    @("aurora",        "Colorado",   14138696,   36.59802652695664)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $row++
}
